$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to reflect the new "through" date
$ws.Name = "Through 2022-08-09"

# Update the column header text for the current month-to-date column
$ws.Range("B1").Value = "August 2022 (through August 09)"

# Apply per-cell carjacking count updates (new data for 2022-08-17 backfilled
# across the matching "August" columns for each prior year)
$ws.Range("AX2").Value = 1
$ws.Range("R3").Value = 1
$ws.Range("R5").Value = 4
$ws.Range("BF6").Value = 2
$ws.Range("R7").Value = 2
$ws.Range("R8").Value = 3
$ws.Range("AP9").Value = 2
$ws.Range("J12").Value = 2
$ws.Range("AX13").Value = 1
$ws.Range("B15").Value = 5
$ws.Range("AH15").Value = 2
$ws.Range("AX17").Value = 3
$ws.Range("B20").Value = 1
$ws.Range("J24").Value = 1
$ws.Range("AX24").Value = 1
$ws.Range("B25").Value = 1
$ws.Range("J27").Value = 1
$ws.Range("J28").Value = 1
$ws.Range("B33").Value = 1
$ws.Range("J37").Value = 1
$ws.Range("B39").Value = 1
$ws.Range("R45").Value = 2
$ws.Range("R46").Value = 1
$ws.Range("AX46").Value = 1
$ws.Range("J64").Value = 2
$ws.Range("J65").Value = 2
$ws.Range("AX69").Value = 1
$ws.Range("R89").Value = 2
$ws.Range("B92").Value = 1
